$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.029099138306892
$ws.Cells.Item(2, 4).Value = 1.032522228300139
$ws.Cells.Item(2, 5).Value = 1.038514979376432
$ws.Cells.Item(2, 6).Value = 1.048498029826579
$ws.Cells.Item(2, 9).Value = 1.031302813553547
$ws.Cells.Item(2, 10).Value = 1.034247791365725
$ws.Cells.Item(2, 11).Value = 1.035327216388917
$ws.Cells.Item(2, 12).Value = 1.041302793615934
$ws.Cells.Item(2, 13).Value = 1.051257707807088
$ws.Cells.Item(2, 14).Value = 1.035716541859344
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.030149717142103
$ws.Cells.Item(3, 4).Value = 1.033305970828889
$ws.Cells.Item(3, 5).Value = 1.039520946324919
$ws.Cells.Item(3, 6).Value = 1.049757326959964
$ws.Cells.Item(3, 9).Value = 1.03149687763318
$ws.Cells.Item(3, 10).Value = 1.034938741569504
$ws.Cells.Item(3, 11).Value = 1.035919815668164
$ws.Cells.Item(3, 12).Value = 1.042118272161966
$ws.Cells.Item(3, 13).Value = 1.052327905474503
$ws.Cells.Item(3, 14).Value = 1.036408473291665
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.030829211760513
$ws.Cells.Item(4, 4).Value = 1.033812274495261
$ws.Cells.Item(4, 5).Value = 1.040172052641252
$ws.Cells.Item(4, 6).Value = 1.050572819888601
$ws.Cells.Item(4, 9).Value = 1.031620303395728
$ws.Cells.Item(4, 10).Value = 1.035384979229168
$ws.Cells.Item(4, 11).Value = 1.036301784345589
$ws.Cells.Item(4, 12).Value = 1.042645512230005
$ws.Cells.Item(4, 13).Value = 1.05302048156009
$ws.Cells.Item(4, 14).Value = 1.03685534465999
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.03111479987688
$ws.Cells.Item(5, 4).Value = 1.034024925280724
$ws.Cells.Item(5, 5).Value = 1.040445820552991
$ws.Cells.Item(5, 6).Value = 1.050915807873124
$ws.Cells.Item(5, 9).Value = 1.03167167729201
$ws.Cells.Item(5, 10).Value = 1.035572373202085
$ws.Cells.Item(5, 11).Value = 1.036462008671425
$ws.Cells.Item(5, 12).Value = 1.042867061354044
$ws.Cells.Item(5, 13).Value = 1.053311661772127
$ws.Cells.Item(5, 14).Value = 1.037043004753843
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.031162747202559
$ws.Cells.Item(6, 4).Value = 1.034060618581891
$ws.Cells.Item(6, 5).Value = 1.040491789903556
$ws.Cells.Item(6, 6).Value = 1.050973406160519
$ws.Cells.Item(6, 9).Value = 1.031680273025672
$ws.Cells.Item(6, 10).Value = 1.035603825490516
$ws.Cells.Item(6, 11).Value = 1.03648889019506
$ws.Cells.Item(6, 12).Value = 1.042904254413734
$ws.Cells.Item(6, 13).Value = 1.053360553469714
$ws.Cells.Item(6, 14).Value = 1.037074501708132
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.030833028083013
$ws.Cells.Item(7, 4).Value = 1.033815116729188
$ws.Cells.Item(7, 5).Value = 1.040175710573142
$ws.Cells.Item(7, 6).Value = 1.050577402301871
$ws.Cells.Item(7, 9).Value = 1.031620991877565
$ws.Cells.Item(7, 10).Value = 1.035387483999607
$ws.Cells.Item(7, 11).Value = 1.036303926667454
$ws.Cells.Item(7, 12).Value = 1.042648472983259
$ws.Cells.Item(7, 13).Value = 1.053024372240339
$ws.Cells.Item(7, 14).Value = 1.036857852987491
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.029454248751214
$ws.Cells.Item(8, 4).Value = 1.032787269329502
$ws.Cells.Item(8, 5).Value = 1.038854913732572
$ws.Cells.Item(8, 6).Value = 1.048923482380027
$ws.Cells.Item(8, 9).Value = 1.03136884267751
$ws.Cells.Item(8, 10).Value = 1.034481477880416
$ws.Cells.Item(8, 11).Value = 1.035527794883016
$ws.Cells.Item(8, 12).Value = 1.041578477328381
$ws.Cells.Item(8, 13).Value = 1.051619368475482
$ws.Cells.Item(8, 14).Value = 1.03595056023569
$ws.Cells.Item(9, 2).Value = 1.019999999999999
$ws.Cells.Item(9, 3).Value = 1.027022344475387
$ws.Cells.Item(9, 4).Value = 1.030969735411825
$ws.Cells.Item(9, 5).Value = 1.036528862216826
$ws.Cells.Item(9, 6).Value = 1.046013960063345
$ws.Cells.Item(9, 9).Value = 1.030908098557654
$ws.Cells.Item(9, 10).Value = 1.032878445270181
$ws.Cells.Item(9, 11).Value = 1.034148810351929
$ws.Cells.Item(9, 12).Value = 1.039689715742605
$ws.Cells.Item(9, 13).Value = 1.049144213364866
$ws.Cells.Item(9, 14).Value = 1.034345251135279
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.025399483986994
$ws.Cells.Item(10, 4).Value = 1.029753813660194
$ws.Cells.Item(10, 5).Value = 1.034979065487839
$ws.Cells.Item(10, 6).Value = 1.044077521638729
$ws.Cells.Item(10, 9).Value = 1.030589912279041
$ws.Cells.Item(10, 10).Value = 1.0318053624192
$ws.Cells.Item(10, 11).Value = 1.033221878577286
$ws.Cells.Item(10, 12).Value = 1.03842831571102
$ws.Cells.Item(10, 13).Value = 1.047494506271031
$ws.Cells.Item(10, 14).Value = 1.033270644383565
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.024696381929572
$ws.Cells.Item(11, 4).Value = 1.029226306926854
$ws.Cells.Item(11, 5).Value = 1.03430819735018
$ws.Cells.Item(11, 6).Value = 1.043239779791099
$ws.Cells.Item(11, 9).Value = 1.030449523055504
$ws.Cells.Item(11, 10).Value = 1.031339662103492
$ws.Cells.Item(11, 11).Value = 1.032818704681343
$ws.Cells.Item(11, 12).Value = 1.037881584051371
$ws.Cells.Item(11, 13).Value = 1.046780250155465
$ws.Cells.Item(11, 14).Value = 1.032804282719992
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.024435158680935
$ws.Cells.Item(12, 4).Value = 1.029030216769983
$ws.Cells.Item(12, 5).Value = 1.034059037292202
$ws.Cells.Item(12, 6).Value = 1.042928716527456
$ws.Cells.Item(12, 9).Value = 1.03039698395995
$ws.Cells.Item(12, 10).Value = 1.031166522540828
$ws.Cells.Item(12, 11).Value = 1.032668676647089
$ws.Cells.Item(12, 12).Value = 1.037678422450307
$ws.Cells.Item(12, 13).Value = 1.046514954475752
$ws.Cells.Item(12, 14).Value = 1.03263089727929
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024491194675779
$ws.Cells.Item(13, 4).Value = 1.02907228558833
$ws.Cells.Item(13, 5).Value = 1.0341124816089
$ws.Cells.Item(13, 6).Value = 1.042995435632488
$ws.Cells.Item(13, 9).Value = 1.030408271522891
$ws.Cells.Item(13, 10).Value = 1.031203668706089
$ws.Cells.Item(13, 11).Value = 1.032700870445341
$ws.Cells.Item(13, 12).Value = 1.037722004981071
$ws.Cells.Item(13, 13).Value = 1.046571860822879
$ws.Cells.Item(13, 14).Value = 1.032668096196366
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.024674790359881
$ws.Cells.Item(14, 4).Value = 1.029210101125461
$ws.Cells.Item(14, 5).Value = 1.034287601087658
$ws.Cells.Item(14, 6).Value = 1.04321406493958
$ws.Cells.Item(14, 9).Value = 1.030445188158852
$ws.Cells.Item(14, 10).Value = 1.031325353552213
$ws.Cells.Item(14, 11).Value = 1.032806308851197
$ws.Cells.Item(14, 12).Value = 1.037864792307331
$ws.Cells.Item(14, 13).Value = 1.046758320515638
$ws.Cells.Item(14, 14).Value = 1.032789953848928
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.024787901766047
$ws.Cells.Item(15, 4).Value = 1.029294993868827
$ws.Cells.Item(15, 5).Value = 1.034395501972414
$ws.Cells.Item(15, 6).Value = 1.043348784393473
$ws.Cells.Item(15, 9).Value = 1.030467881740845
$ws.Cells.Item(15, 10).Value = 1.031400306694351
$ws.Cells.Item(15, 11).Value = 1.032871236985571
$ws.Cells.Item(15, 12).Value = 1.037952757529742
$ws.Cells.Item(15, 13).Value = 1.04687320589721
$ws.Cells.Item(15, 14).Value = 1.032865013433126
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.025446138209183
$ws.Cells.Item(16, 4).Value = 1.029788801401664
$ws.Cells.Item(16, 5).Value = 1.035023593060503
$ws.Cells.Item(16, 6).Value = 1.044133135485055
$ws.Cells.Item(16, 9).Value = 1.03059917442561
$ws.Cells.Item(16, 10).Value = 1.031836247308576
$ws.Cells.Item(16, 11).Value = 1.033248597851042
$ws.Cells.Item(16, 12).Value = 1.038464589138162
$ws.Cells.Item(16, 13).Value = 1.047541910665432
$ws.Cells.Item(16, 14).Value = 1.033301573133027
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.025858926803113
$ws.Cells.Item(17, 4).Value = 1.03009828505093
$ws.Cells.Item(17, 5).Value = 1.035417632347701
$ws.Cells.Item(17, 6).Value = 1.044625337870054
$ws.Cells.Item(17, 9).Value = 1.030680831575416
$ws.Cells.Item(17, 10).Value = 1.032109420506196
$ws.Cells.Item(17, 11).Value = 1.033484822810185
$ws.Cells.Item(17, 12).Value = 1.038785503484753
$ws.Cells.Item(17, 13).Value = 1.047961391757775
$ws.Cells.Item(17, 14).Value = 1.03357513426792
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.026099661567016
$ws.Cells.Item(18, 4).Value = 1.030278704682378
$ws.Cells.Item(18, 5).Value = 1.035647488516051
$ws.Cells.Item(18, 6).Value = 1.044912503515593
$ws.Cells.Item(18, 9).Value = 1.030728208731011
$ws.Cells.Item(18, 10).Value = 1.032268656674731
$ws.Cells.Item(18, 11).Value = 1.033622434543695
$ws.Cells.Item(18, 12).Value = 1.038972635597947
$ws.Cells.Item(18, 13).Value = 1.048206075560622
$ws.Cells.Item(18, 14).Value = 1.033734596570079
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.026181739507376
$ws.Cells.Item(19, 4).Value = 1.030340206662568
$ws.Cells.Item(19, 5).Value = 1.035725866886777
$ws.Cells.Item(19, 6).Value = 1.045010431881166
$ws.Cells.Item(19, 9).Value = 1.030744320353468
$ws.Cells.Item(19, 10).Value = 1.03232293494551
$ws.Cells.Item(19, 11).Value = 1.033669327023607
$ws.Cells.Item(19, 12).Value = 1.039036434039928
$ws.Cells.Item(19, 13).Value = 1.048289507725942
$ws.Cells.Item(19, 14).Value = 1.033788951922228
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.025814642413244
$ws.Cells.Item(20, 4).Value = 1.030065090409742
$ws.Cells.Item(20, 5).Value = 1.035375353624594
$ws.Cells.Item(20, 6).Value = 1.044572521739683
$ws.Cells.Item(20, 9).Value = 1.030672096608104
$ws.Cells.Item(20, 10).Value = 1.032080122069816
$ws.Cells.Item(20, 11).Value = 1.033459496134423
$ws.Cells.Item(20, 12).Value = 1.038751077772832
$ws.Cells.Item(20, 13).Value = 1.047916384628638
$ws.Cells.Item(20, 14).Value = 1.0335457942244
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024620727665476
$ws.Cells.Item(21, 4).Value = 1.029169522041916
$ws.Cells.Item(21, 5).Value = 1.034236031934579
$ws.Cells.Item(21, 6).Value = 1.043149680992921
$ws.Cells.Item(21, 9).Value = 1.030434327957114
$ws.Cells.Item(21, 10).Value = 1.031289524753922
$ws.Cells.Item(21, 11).Value = 1.032775267357212
$ws.Cells.Item(21, 12).Value = 1.037822747241351
$ws.Cells.Item(21, 13).Value = 1.046703412491145
$ws.Cells.Item(21, 14).Value = 1.032754074169634
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.023869719450958
$ws.Cells.Item(22, 4).Value = 1.028605571059504
$ws.Cells.Item(22, 5).Value = 1.033519869907224
$ws.Cells.Item(22, 6).Value = 1.042255728266409
$ws.Cells.Item(22, 9).Value = 1.03028256381301
$ws.Cells.Item(22, 10).Value = 1.030791532147268
$ws.Cells.Item(22, 11).Value = 1.032343495953494
$ws.Cells.Item(22, 12).Value = 1.037238599444697
$ws.Cells.Item(22, 13).Value = 1.045940831072057
$ws.Cells.Item(22, 14).Value = 1.032255374356359
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.024267876250916
$ws.Cells.Item(23, 4).Value = 1.028904614680475
$ws.Cells.Item(23, 5).Value = 1.033899504527447
$ws.Cells.Item(23, 6).Value = 1.042729568883861
$ws.Cells.Item(23, 9).Value = 1.03036323191839
$ws.Cells.Item(23, 10).Value = 1.031055614041033
$ws.Cells.Item(23, 11).Value = 1.032572534901284
$ws.Cells.Item(23, 12).Value = 1.037548311894756
$ws.Cells.Item(23, 13).Value = 1.046345084295711
$ws.Cells.Item(23, 14).Value = 1.032519831276704
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.025834652739414
$ws.Cells.Item(24, 4).Value = 1.030080089938668
$ws.Cells.Item(24, 5).Value = 1.035394457498533
$ws.Cells.Item(24, 6).Value = 1.0445963868517
$ws.Cells.Item(24, 9).Value = 1.030676044342654
$ws.Cells.Item(24, 10).Value = 1.032093361084384
$ws.Cells.Item(24, 11).Value = 1.033470940707405
$ws.Cells.Item(24, 12).Value = 1.038766633430023
$ws.Cells.Item(24, 13).Value = 1.047936721389237
$ws.Cells.Item(24, 14).Value = 1.033559052039887
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.027651328077738
$ws.Cells.Item(25, 4).Value = 1.031440359431726
$ws.Cells.Item(25, 5).Value = 1.03713004207031
$ws.Cells.Item(25, 6).Value = 1.046765566118119
$ws.Cells.Item(25, 9).Value = 1.03102915633779
$ws.Cells.Item(25, 10).Value = 1.033293642044724
$ws.Cells.Item(25, 11).Value = 1.034506653155251
$ws.Cells.Item(25, 12).Value = 1.040178397268065
$ws.Cells.Item(25, 13).Value = 1.049784027872263
$ws.Cells.Item(25, 14).Value = 1.034761037536865
